$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from teste_7gb to teste_7gb_2
$ws.Name = "teste_7gb_2"

# Update the raw measurement data in A2:D51 with the new benchmark run values.
# F2:I2 hold shared AVERAGE() formulas and will recalculate automatically.
$data = New-Object "object[,]" 50,4
$data[0,0] = 48.32
$data[0,1] = 0.02
$data[0,2] = 11.56
$data[0,3] = 0.23
$data[1,0] = 28.37
$data[1,1] = 0.02
$data[1,2] = 10.83
$data[1,3] = 0.38
$data[2,0] = 31.64
$data[2,1] = 0.04
$data[2,2] = 11.66
$data[2,3] = 0.37
$data[3,0] = 28.19
$data[3,1] = 0.03
$data[3,2] = 11.11
$data[3,3] = 0.39
$data[4,0] = 25.99
$data[4,1] = 0.03
$data[4,2] = 10.75
$data[4,3] = 0.41
$data[5,0] = 27.54
$data[5,1] = 0.03
$data[5,2] = 10.85
$data[5,3] = 0.39
$data[6,0] = 27.59
$data[6,1] = 0.02
$data[6,2] = 12.6
$data[6,3] = 0.45
$data[7,0] = 26.3
$data[7,1] = 0.03
$data[7,2] = 10.86
$data[7,3] = 0.41
$data[8,0] = 26.3
$data[8,1] = 0.04
$data[8,2] = 11.49
$data[8,3] = 0.43
$data[9,0] = 25.84
$data[9,1] = 0.03
$data[9,2] = 11.8
$data[9,3] = 0.45
$data[10,0] = 25.75
$data[10,1] = 0.03
$data[10,2] = 11.08
$data[10,3] = 0.43
$data[11,0] = 25.79
$data[11,1] = 0.03
$data[11,2] = 11.86
$data[11,3] = 0.46
$data[12,0] = 25.2
$data[12,1] = 0.02
$data[12,2] = 11.08
$data[12,3] = 0.44
$data[13,0] = 23.96
$data[13,1] = 0.02
$data[13,2] = 10.45
$data[13,3] = 0.43
$data[14,0] = 24.86
$data[14,1] = 0.02
$data[14,2] = 10.84
$data[14,3] = 0.43
$data[15,0] = 23.96
$data[15,1] = 0.03
$data[15,2] = 11.11
$data[15,3] = 0.46
$data[16,0] = 24.76
$data[16,1] = 0.03
$data[16,2] = 11.97
$data[16,3] = 0.48
$data[17,0] = 24.16
$data[17,1] = 0.03
$data[17,2] = 11.13
$data[17,3] = 0.46
$data[18,0] = 24.26
$data[18,1] = 0.03
$data[18,2] = 11.58
$data[18,3] = 0.47
$data[19,0] = 23.49
$data[19,1] = 0.03
$data[19,2] = 10.61
$data[19,3] = 0.45
$data[20,0] = 23.56
$data[20,1] = 0.02
$data[20,2] = 10.46
$data[20,3] = 0.44
$data[21,0] = 23.3
$data[21,1] = 0.02
$data[21,2] = 10.93
$data[21,3] = 0.47
$data[22,0] = 23.94
$data[22,1] = 0.03
$data[22,2] = 11.54
$data[22,3] = 0.48
$data[23,0] = 23.14
$data[23,1] = 0.02
$data[23,2] = 10.68
$data[23,3] = 0.46
$data[24,0] = 23.05
$data[24,1] = 0.02
$data[24,2] = 10.46
$data[24,3] = 0.45
$data[25,0] = 22.24
$data[25,1] = 0.02
$data[25,2] = 10.49
$data[25,3] = 0.47
$data[26,0] = 22.36
$data[26,1] = 0.03
$data[26,2] = 10.35
$data[26,3] = 0.46
$data[27,0] = 22.51
$data[27,1] = 0.02
$data[27,2] = 10.22
$data[27,3] = 0.45
$data[28,0] = 22.25
$data[28,1] = 0.03
$data[28,2] = 10.74
$data[28,3] = 0.48
$data[29,0] = 22.73
$data[29,1] = 0.03
$data[29,2] = 11.16
$data[29,3] = 0.49
$data[30,0] = 22.09
$data[30,1] = 0.02
$data[30,2] = 10.73
$data[30,3] = 0.48
$data[31,0] = 21.35
$data[31,1] = 0.04
$data[31,2] = 10.33
$data[31,3] = 0.48
$data[32,0] = 21.38
$data[32,1] = 0.02
$data[32,2] = 10.21
$data[32,3] = 0.47
$data[33,0] = 21.26
$data[33,1] = 0.03
$data[33,2] = 10.27
$data[33,3] = 0.48
$data[34,0] = 21.28
$data[34,1] = 0.03
$data[34,2] = 10.02
$data[34,3] = 0.47
$data[35,0] = 21.41
$data[35,1] = 0.03
$data[35,2] = 10.4
$data[35,3] = 0.48
$data[36,0] = 21.75
$data[36,1] = 0.03
$data[36,2] = 10.64
$data[36,3] = 0.49
$data[37,0] = 22.44
$data[37,1] = 0.02
$data[37,2] = 11.4
$data[37,3] = 0.5
$data[38,0] = 21.72
$data[38,1] = 0.02
$data[38,2] = 10.41
$data[38,3] = 0.48
$data[39,0] = 20.92
$data[39,1] = 0.02
$data[39,2] = 10.04
$data[39,3] = 0.48
$data[40,0] = 21.31
$data[40,1] = 0.02
$data[40,2] = 10.27
$data[40,3] = 0.48
$data[41,0] = 20.97
$data[41,1] = 0.02
$data[41,2] = 10.12
$data[41,3] = 0.48
$data[42,0] = 20.75
$data[42,1] = 0.03
$data[42,2] = 10.02
$data[42,3] = 0.48
$data[43,0] = 20.97
$data[43,1] = 0.03
$data[43,2] = 10.3
$data[43,3] = 0.49
$data[44,0] = 21.42
$data[44,1] = 0.03
$data[44,2] = 10.96
$data[44,3] = 0.51
$data[45,0] = 21.85
$data[45,1] = 0.03
$data[45,2] = 11.5
$data[45,3] = 0.52
$data[46,0] = 20.95
$data[46,1] = 0.04
$data[46,2] = 10.84
$data[46,3] = 0.51
$data[47,0] = 21.29
$data[47,1] = 0.02
$data[47,2] = 11.05
$data[47,3] = 0.51
$data[48,0] = 20.96
$data[48,1] = 0.02
$data[48,2] = 10.65
$data[48,3] = 0.5
$data[49,0] = 20.73
$data[49,1] = 0.01
$data[49,2] = 10.23
$data[49,3] = 0.49

$ws.Range("A2:D51").Value = $data

# Force recalculation so the shared AVERAGE formulas in F2:I2 refresh to match the new data
$excel.Calculate()

# Refresh the F:I column widths (Excel auto-adjusted these after the new values changed
# the rendered digit counts). Values below are the closest settings the ColumnWidth
# property (quantized to the workbook default font grid) can produce to the saved widths.
$ws.Range("F1").ColumnWidth = 16
$ws.Range("G1").ColumnWidth = 16.857143
$ws.Range("H1").ColumnWidth = 17.285714
$ws.Range("I1").ColumnWidth = 22
